# Actualizacion desde MV -datos-
# Adds daily "Tasas de deposito en bolsa" rows for 03-08-2021 .. 02-09-2021
# (rows 147-169) to Sheet1, mirroring the structure of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 147; Date = "03-08-2021"; B = 0.63; C = 0.96; D = 1.56; E = 1.92; F = -1.26 }
    @{ Row = 148; Date = "04-08-2021"; B = 0.58; C = 0.77; D = 1.56; E = 1.92; F = -1.24 }
    @{ Row = 149; Date = "05-08-2021"; B = 0.72; C = 0.77; D = 1.56; E = 2.04; F = -1.18 }
    @{ Row = 150; Date = "06-08-2021"; B = 0.85; C = 1.04; D = 1.56; E = 2.04; F = -1.57 }
    @{ Row = 151; Date = "09-08-2021"; B = 0.75; C = 1.2; D = 1.52; E = 2.04; F = -1.51 }
    @{ Row = 152; Date = "10-08-2021"; B = 0.72; C = 1.08; D = 1.68; E = 2.04; F = -0.44 }
    @{ Row = 153; Date = "11-08-2021"; B = 0.71; C = 1.08; D = 1.73; E = 2.04; F = -1.41 }
    @{ Row = 154; Date = "12-08-2021"; B = 0.73; C = 1.16; D = 1.71; E = 2.04; F = -0.98 }
    @{ Row = 155; Date = "13-08-2021"; B = 0.73; C = 1.08; D = 1.76; E = 2.04; F = -1.06 }
    @{ Row = 156; Date = "16-08-2021"; B = 0.71; C = 1.31; D = 1.68; E = 2.16; F = -1.2 }
    @{ Row = 157; Date = "17-08-2021"; B = 0.75; C = 1.2; D = 1.68; E = 2.37; F = -1.12 }
    @{ Row = 158; Date = "18-08-2021"; B = 0.82; C = 1.36; D = 1.68; E = 2.37; F = -1.13 }
    @{ Row = 159; Date = "19-08-2021"; B = 0.75; C = 1.37; D = 1.83; E = 2.28; F = -1.09 }
    @{ Row = 160; Date = "20-08-2021"; B = 0.76; C = 1.37; D = 1.68; E = 2.28; F = -1.1 }
    @{ Row = 161; Date = "23-08-2021"; B = 0.78; C = 1.37; D = 1.68; E = 2.4; F = -1.09 }
    @{ Row = 162; Date = "24-08-2021"; B = 0.88; C = 1.37; D = 1.68; E = 2.4; F = -1 }
    @{ Row = 163; Date = "25-08-2021"; B = 0.86; C = 1.4; D = 1.68; E = 2.4; F = -0.98 }
    @{ Row = 164; Date = "26-08-2021"; B = 0.88; C = 1.44; D = 1.68; E = 2.4; F = -1.04 }
    @{ Row = 165; Date = "27-08-2021"; B = 0.6899999999999999; C = 1.44; D = 1.68; E = 2.4; F = -0.98 }
    @{ Row = 166; Date = "30-08-2021"; B = 0.96; C = 1.32; D = 1.68; E = 2.4; F = -0.92 }
    @{ Row = 167; Date = "31-08-2021"; B = 0.99; C = 1.2; D = 1.8; E = 2.4; F = -0.87 }
    @{ Row = 168; Date = "01-09-2021"; B = 1.7; C = 2.27; D = 1.8; E = 2.4; F = -1.25 }
    @{ Row = 169; Date = "02-09-2021"; B = 1.45; C = 2.16; D = 1.8; E = 2.4; F = -0.91 }
)

$firstRow = $rows[0].Row
$lastRow = $rows[$rows.Count - 1].Row

foreach ($r in $rows) {
    # Write the date as a formula that evaluates to a text string, so Excel
    # stores it as literal text ("dd-mm-yyyy") instead of auto-converting it
    # into a date serial number. This mirrors the column A values already
    # present in the sheet (e.g. "02-08-2021").
    $ws.Cells.Item($r.Row, 1).Formula = '="' + $r.Date + '"'

    $ws.Cells.Item($r.Row, 2).Value2 = $r.B
    $ws.Cells.Item($r.Row, 3).Value2 = $r.C
    $ws.Cells.Item($r.Row, 4).Value2 = $r.D
    $ws.Cells.Item($r.Row, 5).Value2 = $r.E
    $ws.Cells.Item($r.Row, 6).Value2 = $r.F
}

# Convert the formula-based text cells in column A into plain shared-string
# values (remove the formula, keep the displayed text), same cell type as
# the other "Serie" cells in the column.
$dateRange = $ws.Range("A" + $firstRow + ":A" + $lastRow)
$dateRange.Copy()
$dateRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$excel.CutCopyMode = 0
